# "9th Stab - Cosmetic Changes"
# Two new rank columns (Jun_17, Jun_15) are inserted right before the
# existing rank column, pushing the old "Jun_13"/"Jun_10" columns two
# slots to the right. The new columns are filled with the same "UN"
# marker used throughout the existing rank column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column B.
# Old column B ("Jun_13" header + data) shifts to D.
# Old column C ("Jun_10" header + data, incl. the highlighted price-target
# cell on row 18) shifts to E.
$ws.Range("B1:C1").EntireColumn.Insert()

# Header labels for the two freshly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Data rows: same "UN" marker as the rest of the rank columns.
$ws.Range("B2:C27").Value = "UN"

# Give the new columns the same width as the pre-existing rank column.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.14
